$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 / Q1, continuing the 0..15 sequence ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header style (bold font + border + center/top alignment) from O1
# onto the two newly added header cells so they pick up style index s="1".
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# --- Data rows 2-25 ---
# For every data row: I flips 1->2, K flips 2->1, M flips 1->2, O flips 2->1,
# and two new trailing columns P and Q are appended, both valued 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
